{"js": "// Replace the multiplication expressions in the table cells with their\n// new values, per the commit's diff. Each value is unique in the\n// document, so a literal, case-sensitive search-and-replace on the whole\n// body safely targets exactly one run each. The pairs are applied in\n// document order, which also guarantees that the one new value that\n// collides with an original value elsewhere (\"65\u00d731=\") is written only\n// after the original occurrence of that same text has already been\n// replaced.\nconst replacements = [\n  [\"83\u00d793=\", \"12\u00d775=\"],\n  [\"47\u00d791=\", \"58\u00d786=\"],\n  [\"44\u00d715=\", \"13\u00d755=\"],\n  [\"72\u00d780=\", \"95\u00d757=\"],\n  [\"74\u00d742=\", \"93\u00d731=\"],\n  [\"16\u00d718=\", \"54\u00d750=\"],\n  [\"18\u00d726=\", \"48\u00d723=\"],\n  [\"51\u00d742=\", \"19\u00d786=\"],\n  [\"52\u00d794=\", \"31\u00d711=\"],\n  [\"54\u00d763=\", \"79\u00d734=\"],\n  [\"65\u00d731=\", \"17\u00d724=\"],\n  [\"14\u00d782=\", \"55\u00d769=\"],\n  [\"91\u00d734=\", \"29\u00d742=\"],\n  [\"23\u00d734=\", \"47\u00d756=\"],\n  [\"29\u00d784=\", \"78\u00d769=\"],\n  [\"26\u00d780=\", \"32\u00d715=\"],\n  [\"33\u00d777=\", \"18\u00d788=\"],\n  [\"54\u00d769=\", \"91\u00d758=\"],\n  [\"61\u00d783=\", \"64\u00d765=\"],\n  [\"45\u00d769=\", \"64\u00d745=\"],\n  [\"77\u00d713=\", \"17\u00d711=\"],\n  [\"45\u00d730=\", \"93\u00d733=\"],\n  [\"97\u00d729=\", \"44\u00d793=\"],\n  [\"92\u00d765=\", \"65\u00d731=\"],\n  [\"42\u00d768=\", \"76\u00d789=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication expressions in the table cells with their\n# new values, per the commit's diff. Each value is unique in the\n# document, so a literal, case-sensitive Find/Replace targets exactly\n# one run each. The pairs are applied in document order, which also\n# guarantees that the one new value that collides with an original\n# value elsewhere (\"65\u00d731=\") is written only after the original\n# occurrence of that same text has already been replaced.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"83\u00d793=\", \"12\u00d775=\"),\n    @(\"47\u00d791=\", \"58\u00d786=\"),\n    @(\"44\u00d715=\", \"13\u00d755=\"),\n    @(\"72\u00d780=\", \"95\u00d757=\"),\n    @(\"74\u00d742=\", \"93\u00d731=\"),\n    @(\"16\u00d718=\", \"54\u00d750=\"),\n    @(\"18\u00d726=\", \"48\u00d723=\"),\n    @(\"51\u00d742=\", \"19\u00d786=\"),\n    @(\"52\u00d794=\", \"31\u00d711=\"),\n    @(\"54\u00d763=\", \"79\u00d734=\"),\n    @(\"65\u00d731=\", \"17\u00d724=\"),\n    @(\"14\u00d782=\", \"55\u00d769=\"),\n    @(\"91\u00d734=\", \"29\u00d742=\"),\n    @(\"23\u00d734=\", \"47\u00d756=\"),\n    @(\"29\u00d784=\", \"78\u00d769=\"),\n    @(\"26\u00d780=\", \"32\u00d715=\"),\n    @(\"33\u00d777=\", \"18\u00d788=\"),\n    @(\"54\u00d769=\", \"91\u00d758=\"),\n    @(\"61\u00d783=\", \"64\u00d765=\"),\n    @(\"45\u00d769=\", \"64\u00d745=\"),\n    @(\"77\u00d713=\", \"17\u00d711=\"),\n    @(\"45\u00d730=\", \"93\u00d733=\"),\n    @(\"97\u00d729=\", \"44\u00d793=\"),\n    @(\"92\u00d765=\", \"65\u00d731=\"),\n    @(\"42\u00d768=\", \"76\u00d789=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"Could not find text '$oldText' to replace.\"\n    }\n}\n\n$d.Saved = $false\n"}
